# Update "想去人数" (interest count) values on the "展览" and "全部类型" sheets
# to match the freshly generated data output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - rows 3,4,5 in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1176
$ws1.Range("F4").Value = 2646
$ws1.Range("F5").Value = 232

# Sheet "全部类型" (fourth sheet) - rows 5,6,8 in column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1176
$ws4.Range("F6").Value = 2646
$ws4.Range("F8").Value = 232
